{"js": "// The three exercise items (paragraphs) that must become green (RGB 00B050),\n// identified by their exact paragraph text. These are the 4th, 5th and 6th\n// numbered exercises in the list, currently colored \"auto\" (black), which\n// should be promoted to the same green used by the first three exercises\n// (\"questions that should be asked\").\nconst targetTexts = [\n  \"\u0628\u0627 \u0627\u0633\u062a\u0641\u0627\u062f\u0647 \u0627\u0632 \u0645\u0627\u0698\u0648\u0644 Faker \u0644\u06cc\u0633\u062a\u06cc \u062d\u0627\u0648\u06cc 4 \u062f\u06cc\u06a9\u0634\u0646\u0631\u06cc \u062d\u0627\u0648\u06cc \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u0627\u0634\u062e\u0627\u0635 \u0628\u0633\u0627\u0632\u06cc\u062f.\",\n  \"\u0628\u0631\u0646\u0627\u0645\u0647 \u0627\u06cc \u0628\u0646\u0648\u06cc\u0633\u06cc\u062f \u06a9\u0647 \u062a\u0627\u0631\u06cc\u062e \u0645\u06cc\u0644\u0627\u062f\u06cc \u06a9\u0631\u06cc\u0633\u0645\u0633 \u0631\u0627 \u0628\u0647 \u0634\u0645\u0633\u06cc \u062a\u0628\u062f\u06cc\u0644 \u06a9\u0646\u062f.\",\n  \"\u0628\u0631\u0646\u0627\u0645\u0647 \u0627\u06cc \u0628\u0646\u0648\u06cc\u0633\u06cc\u062f \u06a9\u0647 \u062a\u0646\u0647\u0627 \u0631\u0648\u0632 \u0648 \u0645\u0627\u0647 \u0627\u0632 \u062a\u0627\u0631\u06cc\u062e \u0634\u0645\u0633\u06cc \u0627\u0645\u0631\u0648\u0632 \u0631\u0627 \u0646\u0645\u0627\u06cc\u0634 \u062f\u0647\u062f.\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (targetTexts.includes(text)) {\n    // Setting Font.Color on the paragraph applies it to both the paragraph\n    // mark run properties and every run of text within the paragraph, just\n    // like selecting the whole list item in Word and changing its color.\n    paragraph.font.color = \"#00B050\";\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The three exercise items (paragraphs) that must become green (RGB 00B050),\n# identified by their exact paragraph text. These are the 4th, 5th and 6th\n# numbered exercises in the list, currently colored \"auto\" (black), which\n# should be promoted to the same green already used by the first three\n# exercises (\"questions that should be asked\").\n$targets = @(\n  \"\u0628\u0627 \u0627\u0633\u062a\u0641\u0627\u062f\u0647 \u0627\u0632 \u0645\u0627\u0698\u0648\u0644 Faker \u0644\u06cc\u0633\u062a\u06cc \u062d\u0627\u0648\u06cc 4 \u062f\u06cc\u06a9\u0634\u0646\u0631\u06cc \u062d\u0627\u0648\u06cc \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u0627\u0634\u062e\u0627\u0635 \u0628\u0633\u0627\u0632\u06cc\u062f.\",\n  \"\u0628\u0631\u0646\u0627\u0645\u0647 \u0627\u06cc \u0628\u0646\u0648\u06cc\u0633\u06cc\u062f \u06a9\u0647 \u062a\u0627\u0631\u06cc\u062e \u0645\u06cc\u0644\u0627\u062f\u06cc \u06a9\u0631\u06cc\u0633\u0645\u0633 \u0631\u0627 \u0628\u0647 \u0634\u0645\u0633\u06cc \u062a\u0628\u062f\u06cc\u0644 \u06a9\u0646\u062f.\",\n  \"\u0628\u0631\u0646\u0627\u0645\u0647 \u0627\u06cc \u0628\u0646\u0648\u06cc\u0633\u06cc\u062f \u06a9\u0647 \u062a\u0646\u0647\u0627 \u0631\u0648\u0632 \u0648 \u0645\u0627\u0647 \u0627\u0632 \u062a\u0627\u0631\u06cc\u062e \u0634\u0645\u0633\u06cc \u0627\u0645\u0631\u0648\u0632 \u0631\u0627 \u0646\u0645\u0627\u06cc\u0634 \u062f\u0647\u062f.\"\n)\n\n# Build the OLE (BGR) color value for RGB 00B050 without relying on an RGB()\n# helper: OLE_COLOR packs as 0x00BBGGRR.\n$red = 0x00\n$green = 0xB0\n$blue = 0x50\n$oleGreen = $red -bor ($green * 256) -bor ($blue * 65536)\n\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text.Trim()\n  if ($targets -contains $text) {\n    # Setting Font.Color on the paragraph's Range (which includes the\n    # trailing paragraph mark) recolors both the paragraph-mark run\n    # properties and every run of text within the paragraph, just like\n    # selecting the whole list item in Word and changing its color.\n    $p.Range.Font.Color = $oleGreen\n  }\n}\n"}
